# skicaDB.xlsx — "Dodat rad sa otkazivanjem posla"
# Adds a new "otkazivanjePosla" (job-cancellation) column (J) to the "posao"
# table on Sheet1: header otkazivanjePosla / idPosao / razlog / status,
# styled with a new light-blue fill, plus a matching column width and an
# updated view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column width for column J (10) ---------------------------------
# (target stored width is 21.140625 character-units; the host's pixel
# quantization only exposes discrete steps, so feed the ColumnWidth value
# whose rounded result lands closest to the target)
$ws.Columns.Item(10).ColumnWidth = 20.33

# --- header cell J1: new fill color + text "otkazivanjePosla" -----------
$headerCell = $ws.Cells.Item(1, 10)
$headerCell.Interior.Color = 16764057   # RGB(153,204,255) = FF99CCFF
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108 # xlCenter
$headerCell.VerticalAlignment = -4108   # xlCenter
$headerCell.Value = "otkazivanjePosla"

# --- new data values under the new column --------------------------------
$ws.Cells.Item(2, 10).Value = "idPosao"
$ws.Cells.Item(3, 10).Value = "razlog"
$ws.Cells.Item(4, 10).Value = "status"

# --- update view: scroll so column C is the leftmost visible column and
#     the active selection moves to G11 ------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G11").Select()
